$d = $word.ActiveDocument

function Replace-Text {
    param([string]$Old, [string]$New)
    $range = $d.Content
    $found = $range.Find.Execute($Old, $true, $false, $false, $false, $false, $true, 1, $false, $New, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $Old"
    }
}

# 1. "Advised and Mentored the TabX team ..." - run merge only, text unchanged
Replace-Text "Advised and Mentored the TabX team to ramp up on web development" "Advised and Mentored the TabX team to ramp up on web development"

# 2. "Presented for Tableau's ..." -> "Frequently presented for Tableau's ..."
Replace-Text "Presented for Tableau" "Frequently presented for Tableau"

# 3. "Attended conferences such as CppCon and React Europe, presenting my learnings from each to the company"
Replace-Text "Attended conferences such as CppCon and React Europe, presenting my learnings from each to the company" "Attended conferences such as CppCon and React Europe, presenting learnings to Tabloids"

# 4. "Added a Chrome Extension, Chrome as new Tab to use Tableau as a Chrome Extension"
Replace-Text "Added a Chrome Extension, Chrome as new Tab to use Tableau as a Chrome Extension" "Added a Chrome Extension, Chrome as new Tab, to use Tableau within a Chrome Extension"

# 5. "Data Move to the Web (partnership)" -> "Data Move to the Web Team"
Replace-Text "Data Move to the Web (partnership)" "Data Move to the Web Team"

# 6. "Added the ability to connect to 27 key datasources" - run merge only, text unchanged
Replace-Text "Added the ability to connect to 27 key datasources" "Added the ability to connect to 27 key datasources"

# 7. Move the _GoBack bookmark from the "Mini Whack-A-Moles" paragraph to
#    immediately after "eight figure contracts".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$target = $d.Content
$found = $target.Find.Execute("eight figure contracts", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target.Collapse(0)
    $target.InsertAfter("Z")
    $zStart = $target.Start
    $zRange = $d.Range($zStart, $zStart + 1)
    $d.Bookmarks.Add("_GoBack", $zRange)
    $zRange2 = $d.Range($zStart, $zStart + 1)
    $zRange2.Text = ""
} else {
    Write-Host "NOT FOUND: eight figure contracts"
}

# 8. "First usage of RSA for encryption in our VizClient code" - run merge only
Replace-Text "First usage of RSA for encryption in our VizClient code" "First usage of RSA for encryption in our VizClient code"

# 9. "Detected unsupported datasources and blocks users from interacting in yet-to-be-built areas" - run merge only
Replace-Text "Detected unsupported datasources and blocks users from interacting in yet-to-be-built areas" "Detected unsupported datasources and blocks users from interacting in yet-to-be-built areas"

# 10. "Detected broken datasources, displaying a red bang in those scenarios"
Replace-Text "Detected broken datasources, displaying a red bang in those scenarios" "Partnered to build various features: joins, toggle referential integrity, broken datasource detection, etc."

# 11. "four person" -> "four-person"
Replace-Text "four person" "four-person"

# 12. "Proficient in Java, ..." -> "Proficient in C#, TypeScript, Java, ..."
Replace-Text "Proficient in Java, JavaScript, C++, C, SQL, HTML5, CSS, JSON, XML, JSP" "Proficient in C#, TypeScript, Java, JavaScript, C++, C, SQL, HTML5, CSS, JSON, XML, JSP"
